$d = $word.ActiveDocument

# The "NOMBRE ALUMNO (A)" / "CALIFICACION" roster table (table #2 in the document).
$t = $d.Tables.Item(2)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    # Trim the trailing cell-mark / paragraph-mark characters Word appends to Range.Text.
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $newText
}

# Column 2 = NOMBRE ALUMNO (A), Column 3 = CALIFICACION.
# Row indices below are 1-based COM table rows (row 1 is the header row).

Set-CellText $t 2 2 "XOTLANIHUA ESPINOSA JARED JESUS"
Set-CellText $t 2 3 "5"

Set-CellText $t 3 2 "BRETON VICENTE AMYRA NAHOMY"

Set-CellText $t 4 2 "ZARATE VERGEL PAULINA"
Set-CellText $t 4 3 "SC"

Set-CellText $t 5 2 "FLORES CERVANTES JAVIER"
Set-CellText $t 5 3 "5"

Set-CellText $t 6 2 "JIMENEZ APARICIO YAZMIN"

Set-CellText $t 7 2 "MARROQUIN HERRERA ESMERALDA"

Set-CellText $t 8 2 "MARCIAL MORALES IVAN DE JESUS"

Set-CellText $t 9 2 "TORRES PEREZ CONSTANZA XIMENA"
Set-CellText $t 9 3 "5"

# Remove the last two roster rows (now-duplicated entries) entirely.
$t.Rows.Item(11).Delete()
$t.Rows.Item(10).Delete()
